# Work allocation.xlsx - apply the "shark_news" update:
#  - insert a new "Status" column (E) between Accountable and Note
#  - refresh several Output/Note/Work-involved descriptions
#  - widen column C, add Status column width, mark row 1 dimension/filter growth
#  - add a stray formatted (Comma style) cell at C15
#  - tidy up window view (zoom/selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at E (old D "Accountable" stays, old E "Note" shifts to F)
$ws.Columns.Item(5).Insert()

# 2. Column widths: widen C (no longer "best fit"), size the new Status column E
$ws.Columns.Item(3).ColumnWidth = 72.14
$ws.Columns.Item(5).ColumnWidth = 17.33

# 3. Header row
$ws.Range("A1").Value = 'Output'
$ws.Range("B1").Value = 'Variable'
$ws.Range("C1").Value = 'Work involved'
$ws.Range("D1").Value = 'Accountable'
$ws.Range("E1").Value = 'Status'
$ws.Range("F1").Value = 'Note'

# 4. Row 2 - Bar Chart: Shark attacks by Area
$ws.Range("A2").Value = 'Bar Chart: Shark attacks by Area'
$ws.Range("B2").Value = '1. shark attacks (count)' + [char]10 + '2. area of shark attacks'
$ws.Range("C2").Value = '1. cleaning "Area" column in Shark Attack excel file' + [char]10 + '2. Obtain latitude & Longitude on shark attack' + [char]10 + '3. use groupby python to get count of shark attack in each area' + [char]10 + '4. potential (may need to group the area further? (i.e northern, central, southern WA? - based on longitude)'
$ws.Range("D2").Value = 'Katherine'
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 'Perth & Peel to be considered as metropolitan area' + [char]10 + 'Anything above Yanchep is classified as Northern WA' + [char]10 + 'Anything below Dawsville is Southern WA'

# 5. Row 3 - Bar Chart: Car Accidents by Area
$ws.Range("A3").Value = 'Bar Chart: Car Accidents by Area'
$ws.Range("B3").Value = '1. car accident (count)' + [char]10 + '2. car accident area'
$ws.Range("C3").Value = '1. use latitude & longitude data in the resource file to determine an area' + [char]10 + '2. Possibly: split up area based on northern, central, or southern WA - this will allow us to present a comparable data with shark attacks bar chart' + [char]10 + '3. use groupby python to get count of car accident based on the above 3 areas' + [char]10 + '4. exclude accidents not involving pedestrians' + [char]10 + ''
$ws.Range("D3").Value = 'Mel'
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 'Perth & Peel to be considered as metropolitan area' + [char]10 + 'Anything above Yanchep is classified as Northern WA' + [char]10 + 'Anything below Dawsville is Southern WA'

# 6. Row 4 - Bar Chart/Line Chart: Shark attack by year
$ws.Range("A4").Value = 'Bar Chart/Line Chart: Shark attack by year'
$ws.Range("B4").Value = '1. Shark Attacks counts' + [char]10 + '2. year'
$ws.Range("C4").Value = '1. Filter excel file for 5 years data of shark attack file' + [char]10 + '2. Groupby the data using python to obtain count of shark attacks year on year' + [char]10 + '3. Possibly: group this based on year and also area (i.e last 5 years in northern territory, the trend of shark attack is xxx'
$ws.Range("D4").Value = 'Katherine'
$ws.Range("E4").Value = 1

# 7. Row 5 - Bar Chart/Line Chart: car accident by year
$ws.Range("A5").Value = 'Bar Chart/Line Chart: car accident by year'
$ws.Range("B5").Value = '1. Car Accident counts' + [char]10 + '2. year'
$ws.Range("C5").Value = '1. Filter excel file for 5 years data of car crash file' + [char]10 + '2. Groupby the data using python to obtain count of car accidents year on year' + [char]10 + '3. Possibly: group this based on year and also area (i.e last 5 years in northern territory, the trend of car accident is xxx)'
$ws.Range("D5").Value = 'Mel'
$ws.Range("E5").Value = 1

# 8. Row 6 - Heatmap: Shark Attack (renamed from "Heatmap: Shark")
$ws.Range("A6").Value = 'Heatmap: Shark Attack'
$ws.Range("B6").Value = '1. Locations of attacks in WA' + [char]10 + '2. Fatality of the attack'
$ws.Range("C6").Value = '1. cleaning "Area" column in Shark Attack excel file' + [char]10 + '2. Assigning "Fatality" score (i.e the higher, the deadlier) to assign weighting for the heatmap'
$ws.Range("D6").Value = 'Katherine'
$ws.Range("E6").Value = 1

# 9. Row 7 - Bar Chart: no. of shark stories in news ... (content + source rewritten)
$ws.Range("A7").Value = 'Bar Chart: no. of shark stories in news in a year  vs total of shark attack count '
$ws.Range("B7").Value = '1. No. of days of stories of shark attacks'
$ws.Range("C7").Value = '1. go to 9news.com.au and look for Shark Articles' + [char]10 + '2. Use Parshub to extract all sharks related article into excel spreadsheet' + [char]10 + '3. Do manual cleaning in Excel (adding year Column, manual check on relevancy of the article to shark attack, or whether the article is a repeat topic'
$ws.Range("D7").Value = 'Kevin'
$ws.Range("E7").Value = 1
$ws.Rows.Item(7).RowHeight = 60

# 10. Row 8 - Bar Chart & Heatmap: shark sightings vs shark attack (WA only)
$ws.Range("A8").Value = 'Bar Chart & Heatmap: shark sightings vs shark attack (WA only)'
$ws.Range("B8").Value = '1. no of shark sightings at a certain point in time' + [char]10 + '2. no. of shark attacks in WA at a certain point in time'
$ws.Range("C8").Value = '1. check https://catalogue.data.wa.gov.au to extract data from 2016 to 2020' + [char]10 + '2. Run a python API to extract latitude, longitude, sighting date, etc and save it as a CSV file' + [char]10 + '3. Do a minor cleanup on CSV file via Excel to group up sightings based on year' + [char]10 + '4. In Python, create a bar chart and a heatmap using the clean CSV file'
$ws.Range("D8").Value = 'Kevin'
$ws.Range("E8").Value = 1
$ws.Rows.Item(8).RowHeight = 107.25

# 11. Row 9 - histogram on time of shark attack (whole australia)
$ws.Range("A9").Value = 'histogram on time of shark attack  (whole australia)'
$ws.Range("B9").Value = '1. Grouping time of attack' + [char]10 + '2. count of shark attacks at certain point in time'
$ws.Range("C9").Value = '1. cleaning up the time column in shark attack data' + [char]10 + '2. Grouping the time to Dawn (4am to 8am), Daytime (8am to 4pm)  Dusk (4pm to 8pm), Nightime (8pm to 4am)'
$ws.Range("D9").Value = 'Cathy'
$ws.Range("E9").Value = 1

# 12. Row 10 - bar chart: human activity vs shark attack (whole Australia)
$ws.Range("A10").Value = 'bar chart: human activity vs shark attack (whole Australia)'
$ws.Range("B10").Value = '1. Group of activities being performed (fishing, swimming, diving)' + [char]10 + '2. count of shark attacks in relation to the activity'
$ws.Range("C10").Value = '1. Using the Shark Attacks excel spreadsheet, do a manual cleaning on the activity column by grouping up activities into smaller group' + [char]10 + '2. run a bar chart to see the frequency of shark attacks based on human activities undertaken'
$ws.Range("D10").Value = 'Cathy'
$ws.Range("E10").Value = 1

# 13. Format the new "Status" column (E2:E10) as a bordered percentage cell
$statusRange = $ws.Range("E2:E10")
$statusRange.NumberFormat = "0%"
$statusRange.Borders.LineStyle = 1

# 14. Stray formatted cell below the table (Comma number style)
$ws.Range("C15").Style = "Comma"

# 15. Re-apply the AutoFilter over the new A1:F10 range and fix up the hidden
#     _FilterDatabase defined name that backs it
$ws.AutoFilterMode = $false
$ws.Range("A1:F10").AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$10"
    }
}

# 16. Window view: zoom to 70% and move the selection to E6
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("E6").Select()
